$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "to open (something) (～を)"
$ws.Range("A3").Value = "to close (something) (～を)"
$ws.Range("A4").Value = "to teach; to instruct (person に thing を)"
$ws.Range("A5").Value = "to forget; to leave behind (～を)"
$ws.Range("A6").Value = "to get off (～を)"
$ws.Range("A7").Value = "to borrow (person に thing を)"
$ws.Range("A9").Value = "to turn on (～を)"
$ws.Range("A10").Value = "to call (～に)"
$ws.Range("A11").Value = "to bring (a person) (～を)"
$ws.Range("A12").Value = "to bring (a thing) (～を)"
$ws.Range("A46").Value = "to return (a thing) (person に thing を)"
$ws.Range("A47").Value = "to turn off; to erase (～を)"
$ws.Range("A49").Value = "to sit down (seat に)"
$ws.Range("A52").Value = "to use (～を)"
$ws.Range("A53").Value = "to help (person/task を)"
$ws.Range("A54").Value = "to enter (～に)"
$ws.Range("A55").Value = "to carry; to hold (～を)"
$ws.Range("A56").Value = "(1) to be absent (from...) (～を); (2) to rest"
